# "Generate Report for Archive"
#
# 1) Status text "Ready for handoff" -> "In Translation" everywhere it
#    appears (Overview!E2:F4 = zh-cn/de-de status columns, and the
#    per-language sheets' Status column, zh-cn!C2:C4 / de-de!C2:C4).
# 2) The Status-ish columns got narrower: Overview columns E & F, and
#    column C on the zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: zh-cn (E) / de-de (F) status columns, rows 2-4 ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2:F4").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- zh-cn sheet: Status column (C), rows 2-4 ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2:C4").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- de-de sheet: Status column (C), rows 2-4 ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2:C4").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
